$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting so that
# numeric-looking values (e.g. "1.00") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.838.42'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '2.606.74'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '585.46'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').Value = '164.99'
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.526'
$ws.Range('E8').Value = '  -3.76%  '
$ws.Range('D9').Value = '2.604.65'
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('E10').Value = '  -3.64%  '
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('D14').Value = '27.06'
$ws.Range('E14').Value = '  -3.63%  '
$ws.Range('D15').Value = '3.082.65'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').Value = '0.0000178'
$ws.Range('E16').Value = '  -3.30%  '
$ws.Range('D17').Value = '66.742.11'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '2.620.03'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').Value = '11.62'
$ws.Range('E19').Value = '  -4.43%  '
$ws.Range('D20').Value = '7.77'
$ws.Range('E20').Value = '  -5.10%  '
$ws.Range('D21').Value = '353.28'
$ws.Range('E21').Value = '  -2.76%  '
$ws.Range('E22').Value = '  -3.62%  '
$ws.Range('E23').Value = '  -4.22%  '
$ws.Range('D24').Value = '10.47'
$ws.Range('E24').Value = '  -5.76%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  -5.96%  '
$ws.Range('D27').Value = '69.08'
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('D28').Value = '2.742.84'
$ws.Range('E28').Value = '  -1.64%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = '0.0₃0984'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('D31').Value = '537.16'
$ws.Range('E31').Value = '  -3.81%  '
$ws.Range('D32').Value = '8.11'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('E33').Value = '  -4.63%  '
$ws.Range('E34').Value = '  -3.52%  '
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').Value = '1.46'
$ws.Range('E37').Value = '  -5.44%  '
$ws.Range('D38').Value = '158.01'
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('D39').Value = '18.82'
$ws.Range('E39').Value = '  -2.73%  '
$ws.Range('D40').Value = '0.362'
$ws.Range('E40').Value = '  -2.93%  '
$ws.Range('D41').Value = '18.24'
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('E42').Value = '  -1.89%  '
$ws.Range('D43').Value = '5.09'
$ws.Range('E43').Value = '  -4.51%  '
$ws.Range('E45').Value = '  -6.06%  '
$ws.Range('E46').Value = '  -2.81%  '
$ws.Range('D47').Value = '150.05'
$ws.Range('E47').Value = '  -2.77%  '
$ws.Range('D48').Value = '0.571'
$ws.Range('E48').Value = '  -4.03%  '
$ws.Range('E49').Value = '  -3.67%  '
$ws.Range('D50').Value = '1.69'
$ws.Range('E50').Value = '  -2.40%  '
$ws.Range('E51').Value = '  -1.67%  '
